$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3194513333333333
$ws.Range("H2").Value = 0.9583539999999999
$ws.Range("I2").Value = 0.01095865642710367
$ws.Range("J2").Value = 0.01095865642710367
$ws.Range("M2").Value = 33.54649666666666
$ws.Range("N2").Value = 100.63949
$ws.Range("O2").Value = 0.1978943147725085
$ws.Range("P2").Value = 0.1978943147725085
$ws.Range("Q2").Value = 10.71647308882889
$ws.Range("R2").Value = 96.44825779945999
$ws.Range("S2").Value = 0.002168655804469028
$ws.Range("T2").Value = 0.002168655804469028

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3194513333333333
$ws.Range("H3").Value = 0.9583539999999999
$ws.Range("I3").Value = 0.01095865642710367
$ws.Range("J3").Value = 0.01095865642710367
$ws.Range("O3").Value = 0.6322484766686425
$ws.Range("P3").Value = 0.6322484766686425
$ws.Range("Q3").Value = 34.23783949256644
$ws.Range("R3").Value = 308.140555433098
$ws.Range("S3").Value = 0.006928593832371325
$ws.Range("T3").Value = 0.006928593832371327

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3194513333333333
$ws.Range("H4").Value = 0.9583539999999999
$ws.Range("I4").Value = 0.01095865642710367
$ws.Range("J4").Value = 0.01095865642710367
$ws.Range("O4").Value = 0.169857208558849
$ws.Range("P4").Value = 0.169857208558849
$ws.Range("Q4").Value = 9.198193523431998
$ws.Range("R4").Value = 82.78374171088798
$ws.Range("S4").Value = 0.00186140679026332
$ws.Range("T4").Value = 0.00186140679026332

# Row 5
$ws.Range("I5").Value = 0.9713235907985359
$ws.Range("J5").Value = 0.971323590798536
$ws.Range("M5").Value = 33.54649666666666
$ws.Range("N5").Value = 100.63949
$ws.Range("O5").Value = 0.1978943147725085
$ws.Range("P5").Value = 0.1978943147725085
$ws.Range("Q5").Value = 949.8576025791375
$ws.Range("R5").Value = 8548.718423212238
$ws.Range("S5").Value = 0.1922194164234487
$ws.Range("T5").Value = 0.1922194164234488

# Row 6
$ws.Range("I6").Value = 0.9713235907985359
$ws.Range("J6").Value = 0.971323590798536
$ws.Range("O6").Value = 0.6322484766686425
$ws.Range("P6").Value = 0.6322484766686425
$ws.Range("S6").Value = 0.6141178606346902
$ws.Range("T6").Value = 0.6141178606346902

# Row 7
$ws.Range("I7").Value = 0.9713235907985359
$ws.Range("J7").Value = 0.971323590798536
$ws.Range("O7").Value = 0.169857208558849
$ws.Range("P7").Value = 0.169857208558849
$ws.Range("S7").Value = 0.164986313740397
$ws.Range("T7").Value = 0.164986313740397

# Row 8
$ws.Range("I8").Value = 0.01771775277436037
$ws.Range("J8").Value = 0.01771775277436037
$ws.Range("M8").Value = 33.54649666666666
$ws.Range("N8").Value = 100.63949
$ws.Range("O8").Value = 0.1978943147725085
$ws.Range("P8").Value = 0.1978943147725085
$ws.Range("Q8").Value = 17.32619523789
$ws.Range("R8").Value = 155.93575714101
$ws.Range("S8").Value = 0.003506242544590758
$ws.Range("T8").Value = 0.003506242544590758

# Row 9
$ws.Range("I9").Value = 0.01771775277436037
$ws.Range("J9").Value = 0.01771775277436037
$ws.Range("O9").Value = 0.6322484766686425
$ws.Range("P9").Value = 0.6322484766686425
$ws.Range("S9").Value = 0.01120202220158096
$ws.Range("T9").Value = 0.01120202220158096

# Row 10
$ws.Range("I10").Value = 0.01771775277436037
$ws.Range("J10").Value = 0.01771775277436037
$ws.Range("O10").Value = 0.169857208558849
$ws.Range("P10").Value = 0.169857208558849
$ws.Range("S10").Value = 0.003009488028188656
$ws.Range("T10").Value = 0.003009488028188655

